# Updated cryptos list on Wed Aug  7 14:58:02 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without leaving a
# lingering number format on the cell itself (style is cleared again
# right after the write), so cells that looked like numbers (e.g.
# "1.00", "0.999", "56.170.80") stay as text, exactly like the source data.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "56.170.80"
Set-TextValue $ws.Range("E2") "  +0.34%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.392.66"
Set-TextValue $ws.Range("E3") "  -4.73%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.17%  "

# Row 5
Set-TextValue $ws.Range("D5") "478.54"
Set-TextValue $ws.Range("E5") "  -2.06%  "

# Row 6
Set-TextValue $ws.Range("D6") "147.78"
Set-TextValue $ws.Range("E6") "  +1.94%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  +0.16%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.500"
Set-TextValue $ws.Range("E8") "  -2.61%  "

# Row 9
Set-TextValue $ws.Range("D9") "2.403.51"
Set-TextValue $ws.Range("E9") "  -5.07%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0979"
Set-TextValue $ws.Range("E10") "  -0.29%  "

# Row 11
Set-TextValue $ws.Range("D11") "5.54"
Set-TextValue $ws.Range("E11") "  -1.56%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.326"
Set-TextValue $ws.Range("E12") "  -2.51%  "

# Row 13
Set-TextValue $ws.Range("E13") "  +0.83%  "

# Row 14
Set-TextValue $ws.Range("D14") "2.816.97"
Set-TextValue $ws.Range("E14") "  -4.51%  "

# Row 15
Set-TextValue $ws.Range("D15") "56.400.64"
Set-TextValue $ws.Range("E15") "  +0.69%  "

# Row 16
Set-TextValue $ws.Range("D16") "20.47"
Set-TextValue $ws.Range("E16") "  -3.04%  "

# Row 17
Set-TextValue $ws.Range("E17") "  -2.25%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.400.34"
Set-TextValue $ws.Range("E18") "  -5.04%  "

# Row 19
Set-TextValue $ws.Range("D19") "4.50"
Set-TextValue $ws.Range("E19") "  +1.59%  "

# Row 20
Set-TextValue $ws.Range("D20") "315.58"
Set-TextValue $ws.Range("E20") "  -1.74%  "

# Row 21
Set-TextValue $ws.Range("D21") "9.79"
Set-TextValue $ws.Range("E21") "  -4.43%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.999"
Set-TextValue $ws.Range("E22") "  +0.10%  "

# Row 23
Set-TextValue $ws.Range("D23") "5.71"
Set-TextValue $ws.Range("E23") "  -1.99%  "

# Row 24
Set-TextValue $ws.Range("D24") "56.62"
Set-TextValue $ws.Range("E24") "  -3.26%  "

# Row 25
Set-TextValue $ws.Range("E25") "  +0.37%  "

# Row 26
Set-TextValue $ws.Range("E26") "  -3.70%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.159"
Set-TextValue $ws.Range("E27") "  -4.77%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.507.82"
Set-TextValue $ws.Range("E28") "  -4.26%  "

# Row 29
Set-TextValue $ws.Range("D29") "7.33"
Set-TextValue $ws.Range("E29") "  -2.12%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0₃0777"
Set-TextValue $ws.Range("E30") "  -2.00%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.999"
Set-TextValue $ws.Range("E31") "  -0.03%  "

# Row 32
Set-TextValue $ws.Range("D32") "148.82"
Set-TextValue $ws.Range("E32") "  -0.27%  "

# Row 33
Set-TextValue $ws.Range("D33") "18.00"
Set-TextValue $ws.Range("E33") "  -2.58%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.49"
Set-TextValue $ws.Range("E34") "  -1.11%  "

# Row 35
Set-TextValue $ws.Range("D35") "5.01"
Set-TextValue $ws.Range("E35") "  -4.53%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.12"
Set-TextValue $ws.Range("E36") "  -3.33%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -1.56%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.62"
Set-TextValue $ws.Range("E38") "  -2.15%  "

# Row 39
Set-TextValue $ws.Range("D39") "33.57"
Set-TextValue $ws.Range("E39") "  -1.92%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.36"
Set-TextValue $ws.Range("E40") "  +2.92%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.998"
Set-TextValue $ws.Range("E41") "  +0.38%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.0548"
Set-TextValue $ws.Range("E42") "  -1.63%  "

# Row 43
Set-TextValue $ws.Range("D43") "3.40"
Set-TextValue $ws.Range("E43") "  -4.11%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.589"
Set-TextValue $ws.Range("E44") "  -4.94%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.0949"
Set-TextValue $ws.Range("E45") "  +4.49%  "

# Row 46
Set-TextValue $ws.Range("E46") "  +0.27%  "

# Row 47
Set-TextValue $ws.Range("D47") "254.57"
Set-TextValue $ws.Range("E47") "  -4.16%  "

# Row 48
Set-TextValue $ws.Range("B48") "RenderToken"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D48") "4.61"
Set-TextValue $ws.Range("E48") "  -4.63%  "

# Row 49
Set-TextValue $ws.Range("B49") "VeChain"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D49") "0.0223"
Set-TextValue $ws.Range("E49") "  -1.43%  "

# Row 50
Set-TextValue $ws.Range("D50") "17.27"
Set-TextValue $ws.Range("E50") "  -3.05%  "

# Row 51
Set-TextValue $ws.Range("B51") "Maker"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D51") "1.789.76"
Set-TextValue $ws.Range("E51") "  -8.46%  "

